# Updated cryptos list values (Price and Volume(1h) columns) to match the
# latest scrape. Each target cell is forced to Text format ("@") before the
# new value is written so that Excel does not reinterpret numeric-looking
# strings (e.g. "19.90", "16.00") as floating point numbers and silently
# drop trailing zeros / significant formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.287.72'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.969.05'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -4.25%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '242.93'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.39%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.54'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -9.22%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.374'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.54%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.02'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -4.96%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0798'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +4.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.103'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.858'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -7.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.46'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.07'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -7.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.256.79'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.44'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.83%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.969.80'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.178.78'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.05%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '71.26'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0857'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '236.97'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.71%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.23'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.48%  '
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.53'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.79%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -4.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.84'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +1.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '160.21'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.90'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.56%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.122'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +7.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.120'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -2.04%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.88'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -7.14%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -5.78%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.41'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -7.72%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.29'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.49%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -5.84%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -2.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.10'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +9.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0992'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.01%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.23'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.44%  '
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.27%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.31%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '92.44'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -3.37%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.00'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -6.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.53'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -6.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.339.86'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -6.08%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.151.61'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.97%  '
